# Auto-generated edit script applying the Gungnir_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 479.5
$ws.Range("I38").Value = 222.66667
$ws.Range("J38").Value = 1250
$ws.Range("K38").Value = 668.00001
$ws.Range("L38").Value = 3750
$ws.Range("M38").Value = -296.00001
$ws.Range("N38").Value = -4494

# Row 62
$ws.Range("H62").Value = 9471.048000000001
$ws.Range("I62").Value = 8650.666999999999
$ws.Range("J62").Value = 11522
$ws.Range("K62").Value = 8650.666999999999
$ws.Range("L62").Value = 11522
$ws.Range("M62").Value = -8026.666999999999

# Row 65
$ws.Range("H65").Value = 9471.048000000001
$ws.Range("I65").Value = 8650.666999999999
$ws.Range("J65").Value = 11522
$ws.Range("K65").Value = 43253.335
$ws.Range("L65").Value = 57610
$ws.Range("M65").Value = -40133.335

# Row 74
$ws.Range("H74").Value = 3874.75
$ws.Range("I74").Value = 3899
$ws.Range("J74").Value = 3866.6667
$ws.Range("K74").Value = 3899
$ws.Range("L74").Value = 3866.6667
$ws.Range("M74").Value = -2963
$ws.Range("N74").Value = -5738.6667

# Row 77
$ws.Range("H77").Value = 3874.75
$ws.Range("I77").Value = 3899
$ws.Range("J77").Value = 3866.6667
$ws.Range("K77").Value = 19495
$ws.Range("L77").Value = 19333.3335
$ws.Range("M77").Value = -14815
$ws.Range("N77").Value = -28693.3335

# Row 82
$ws.Range("H82").Value = 907
$ws.Range("I82").Value = 688.4
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 2065.2
$ws.Range("L82").Value = 6000
$ws.Range("M82").Value = -1659.2
$ws.Range("N82").Value = -6812

# Row 85
$ws.Range("H85").Value = 907
$ws.Range("I85").Value = 688.4
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 2065.2
$ws.Range("L85").Value = 6000
$ws.Range("M85").Value = -661.1999999999998
$ws.Range("N85").Value = -8808

# Row 96
$ws.Range("H96").Value = 1603.8
$ws.Range("I96").Value = 996.6667
$ws.Range("J96").Value = 2514.5
$ws.Range("K96").Value = 2990.0001
$ws.Range("L96").Value = 7543.5
$ws.Range("M96").Value = -1617.0001
$ws.Range("N96").Value = -10289.5

# Row 99
$ws.Range("H99").Value = 342.5
$ws.Range("I99").Value = 364.2857
$ws.Range("J99").Value = 190
$ws.Range("K99").Value = 1092.8571
$ws.Range("L99").Value = 570
$ws.Range("M99").Value = 405.1428999999998
$ws.Range("N99").Value = -3566

# Row 101
$ws.Range("H101").Value = 11364128
$ws.Range("I101").Value = 584
$ws.Range("J101").Value = 22727672
$ws.Range("K101").Value = 1752
$ws.Range("L101").Value = 68183016
$ws.Range("M101").Value = -130
$ws.Range("N101").Value = -68186260

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6945842.5
$ws.Range("I61").Value = 8334717
$ws.Range("J61").Value = 1469
$ws.Range("K61").Value = 8334717
$ws.Range("L61").Value = 1469
$ws.Range("M61").Value = -8334505
$ws.Range("N61").Value = -1893

# Row 64
$ws.Range("H64").Value = 17212.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 17212.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 17212.5
$ws.Range("N64").Value = -17708.5

# Row 67
$ws.Range("H67").Value = 17212.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 17212.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 17212.5
$ws.Range("N67").Value = -18928.5

# Row 74
$ws.Range("H74").Value = 1563.6072
$ws.Range("I74").Value = 1379.6818
$ws.Range("J74").Value = 2238
$ws.Range("K74").Value = 1379.6818
$ws.Range("L74").Value = 2238
$ws.Range("M74").Value = -505.6818000000001
$ws.Range("N74").Value = -3986

# Row 77
$ws.Range("H77").Value = 1563.6072
$ws.Range("I77").Value = 1379.6818
$ws.Range("J77").Value = 2238
$ws.Range("K77").Value = 6898.409000000001
$ws.Range("L77").Value = 11190
$ws.Range("M77").Value = -2530.409000000001
$ws.Range("N77").Value = -19926

# Row 122
$ws.Range("H122").Value = 2069.5715
$ws.Range("I122").Value = 2041.375
$ws.Range("J122").Value = 2159.8
$ws.Range("K122").Value = 6124.125
$ws.Range("L122").Value = 6479.400000000001
$ws.Range("M122").Value = -3674.125
$ws.Range("N122").Value = -11379.4

# Row 132
$ws.Range("H132").Value = 2102311.5
$ws.Range("I132").Value = 1182.3478
$ws.Range("J132").Value = 11767505
$ws.Range("K132").Value = 3547.0434
$ws.Range("L132").Value = 35302515
$ws.Range("M132").Value = -1017.0434
$ws.Range("N132").Value = -35307575

# Row 136
$ws.Range("H136").Value = 6945842.5
$ws.Range("I136").Value = 8334717
$ws.Range("J136").Value = 1469
$ws.Range("K136").Value = 25004151
$ws.Range("L136").Value = 4407
$ws.Range("M136").Value = -25001601
$ws.Range("N136").Value = -9507

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5746.1924
$ws.Range("I134").Value = 1880.05
$ws.Range("J134").Value = 18633.334
$ws.Range("K134").Value = 5640.15
$ws.Range("L134").Value = 55900.00199999999
$ws.Range("M134").Value = -3105.15
$ws.Range("N134").Value = -60970.00199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1765.8572
$ws.Range("I16").Value = 772.2
$ws.Range("J16").Value = 4250
$ws.Range("K16").Value = 772.2
$ws.Range("L16").Value = 4250
$ws.Range("M16").Value = -485.2
$ws.Range("N16").Value = -4824

# Row 48
$ws.Range("H48").Value = 7000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 7000
$ws.Range("N48").Value = -7952

# Row 113
$ws.Range("H113").Value = 1765.8572
$ws.Range("I113").Value = 772.2
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 772.2
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 1397.8
$ws.Range("N113").Value = -8590

# Row 134
$ws.Range("H134").Value = 1399.9131
$ws.Range("I134").Value = 1169.4667
$ws.Range("J134").Value = 1832
$ws.Range("K134").Value = 3508.4001
$ws.Range("L134").Value = 5496
$ws.Range("M134").Value = -973.4000999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 50000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 50000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 150000
$ws.Range("N75").Value = -151996
$ws.Range("M75").ClearContents()

# Row 78
$ws.Range("H78").Value = 50000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 50000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 450000
$ws.Range("N78").Value = -459984
$ws.Range("M78").ClearContents()

# Row 113
$ws.Range("H113").Value = 9188528
$ws.Range("I113").Value = 4630157.5
$ws.Range("J113").Value = 19444862
$ws.Range("K113").Value = 13890472.5
$ws.Range("L113").Value = 58334586
$ws.Range("M113").Value = -13888302.5
$ws.Range("N113").Value = -58338926

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

# Row 131
$ws.Range("H131").Value = 806.35
$ws.Range("I131").Value = 625.8
$ws.Range("J131").Value = 815.85266
$ws.Range("K131").Value = 1877.4
$ws.Range("L131").Value = 2447.55798
$ws.Range("M131").Value = 3162.6
$ws.Range("N131").Value = -12527.55798

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 14000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 14000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 14000
$ws.Range("N64").Value = -14496

# Row 67
$ws.Range("H67").Value = 14000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 14000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 14000
$ws.Range("N67").Value = -15716

# Row 100
$ws.Range("H100").Value = 51727.273
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 51727.273
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 51727.273
$ws.Range("N100").Value = -53891.273

# Row 126
$ws.Range("H126").Value = 2534.2727
$ws.Range("I126").Value = 1488
$ws.Range("J126").Value = 2638.9
$ws.Range("K126").Value = 4464
$ws.Range("L126").Value = 7916.700000000001
$ws.Range("M126").Value = -1994
$ws.Range("N126").Value = -12856.7

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4999.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4999.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4999.5
$ws.Range("N7").Value = -5223.5
$ws.Range("M7").ClearContents()

# Row 40
$ws.Range("H40").Value = 11366607
$ws.Range("I40").Value = 3166.6924
$ws.Range("J40").Value = 27780466
$ws.Range("K40").Value = 3166.6924
$ws.Range("L40").Value = 27780466
$ws.Range("M40").Value = -3030.6924
$ws.Range("N40").Value = -27780738

# Row 43
$ws.Range("H43").Value = 336000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 336000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 336000
$ws.Range("N43").Value = -336386

# Row 61
$ws.Range("H61").Value = 2116.818
$ws.Range("I61").Value = 2116.818
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2116.818
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1914.818

# Row 113
$ws.Range("H113").Value = 2116.818
$ws.Range("I113").Value = 2116.818
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2116.818
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 53.18199999999979

# Row 122
$ws.Range("H122").Value = 9535.714
$ws.Range("I122").Value = 16785.715
$ws.Range("J122").Value = 2285.7144
$ws.Range("K122").Value = 50357.145
$ws.Range("L122").Value = 6857.1432
$ws.Range("M122").Value = -47907.145

# Row 126
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 14998.5
$ws.Range("N126").Value = -19938.5
$ws.Range("M126").ClearContents()

# Row 136
$ws.Range("H136").Value = 93170010
$ws.Range("I136").Value = 57145108
$ws.Range("J136").Value = 333336000
$ws.Range("K136").Value = 171435324
$ws.Range("L136").Value = 1000008000
$ws.Range("M136").Value = -171432774
$ws.Range("N136").Value = -1000013100

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 22873
$ws.Range("I70").Value = 8095
$ws.Range("J70").Value = 23928.572
$ws.Range("K70").Value = 8095
$ws.Range("L70").Value = 23928.572
$ws.Range("M70").Value = -7780
$ws.Range("N70").Value = -24558.572

# Row 73
$ws.Range("H73").Value = 22873
$ws.Range("I73").Value = 8095
$ws.Range("J73").Value = 23928.572
$ws.Range("K73").Value = 8095
$ws.Range("L73").Value = 23928.572
$ws.Range("M73").Value = -7003
$ws.Range("N73").Value = -26112.572

# Row 113
$ws.Range("H113").Value = 43480460
$ws.Range("I113").Value = 45456796
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 136370388
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -136368218
$ws.Range("N113").Value = -7340

# Row 122
$ws.Range("H122").Value = 27460.3
$ws.Range("I122").Value = 51041.6
$ws.Range("J122").Value = 3879
$ws.Range("K122").Value = 153124.8
$ws.Range("L122").Value = 11637
$ws.Range("M122").Value = -150674.8
$ws.Range("N122").Value = -16537

# Row 126
$ws.Range("H126").Value = 2999.8572
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -18439.25
